$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.868.53"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "'3.628.46"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'605.79"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "'200.02"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +8.82%  "
$ws.Range("D10").Value = "'0.648"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'53.83"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").Value = "'9.56"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'4.208.68"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "'675.56"
$ws.Range("E15").Value = "  +13.47%  "
$ws.Range("D16").Value = "'70.977.07"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "'12.88"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "'3.629.29"
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("D19").Value = "'19.00"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "'18.47"
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "'104.79"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").Value = "'4.62"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "'10.53"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "'9.82"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").Value = "'34.16"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("D30").Value = "'4.67"
$ws.Range("E30").Value = "  +9.67%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("D32").Value = "'12.19"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").Value = "'63.30"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'3.940.79"
$ws.Range("E35").Value = "  +5.26%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "'0.0₃0865"
$ws.Range("E36").Value = "  +7.41%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'514.86"
$ws.Range("E38").Value = "  +4.83%  "
$ws.Range("D39").Value = "'3.01"
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.55"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'36.44"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("D44").Value = "'0.0459"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("E45").Value = "  +8.93%  "
$ws.Range("E46").Value = "  +6.38%  "
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "'8.63"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'0.000247"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("E51").Value = "  +2.37%  "
